$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.683.91'
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.472.73'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.91'
$ws.Range("E5").Value = '  +1.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.97'
$ws.Range("E6").Value = '  +0.99%  '

$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  +0.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.27'
$ws.Range("E10").Value = '  +2.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("E11").Value = '  +9.22%  '

$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.856.56'
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("E15").Value = '  -1.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.482.56'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("E17").Value = '  +2.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.651.29'
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0950'
$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.15'
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.28'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.30'
$ws.Range("E23").Value = '  +1.05%  '

$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +2.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.68'
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("E28").Value = '  +0.95%  '

$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.05'
$ws.Range("E30").Value = '  +1.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.92'
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("E34").Value = '  +0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0766'
$ws.Range("E35").Value = '  +1.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.55'
$ws.Range("E36").Value = '  +1.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  +5.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.92'
$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("E39").Value = '  +1.77%  '

$ws.Range("E40").Value = '  +0.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  +0.40%  '

$ws.Range("E42").Value = '  +11.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.994.19'
$ws.Range("E43").Value = '  +2.42%  '

$ws.Range("E44").Value = '  +0.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.81'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.99'
$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("E47").Value = '  +4.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.714.56'
$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.38'
$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.83'
$ws.Range("E50").Value = '  +3.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.08'
$ws.Range("E51").Value = '  -0.34%  '
